# Move C erosion from climate to other soil properties
# The filtered view previously constrained group_level1 = "Other Soil Properties"
# AND group_level2 = "Chemical Properties". Drop the second (Chemical Properties)
# filter criterion so the view shows all "Other Soil Properties" rows again.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$filterRange = $ws.AutoFilter.Range
$filterRange.AutoFilter(3)

# Remove the trailing "Labile carbon" / "Natural abundance of 13C" rows (Chemical
# Properties entries that no longer belong on this sheet).
$ws.Rows("211:213").Delete()

# Update the active selection to the last real cell of data, matching where the
# user ended up after the edit.
$ws.Range("F210").Select()
